# Generate Report for Handback
# Adds a new handback record (185fe346-231a-4d9d-9928-1b9795caeec0) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the existing
# rows for the other two handoff files.

$wb = $excel.ActiveWorkbook

$newMd   = "185fe346-231a-4d9d-9928-1b9795caeec0.md"
$newXlfBase = "185fe346-231a-4d9d-9928-1b9795caeec0.871fef4d4bd910215edf6a3c3bd98107dcc9063e"
$status  = "Handed back: in sync with en-US"
$reason  = "Include"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f1a2b3c4d5e6f708192a3b4c5d6e7f809182a3b/e2e/$newMd",
    "",
    "",
    $newMd
)
$wsOverview.Range("A4").Style = $wsOverview.Range("A3").Style

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$newXlfBase.zh-cn.xlf"

$wsZhCn.Range("B4").Value = $status
$wsZhCn.Range("D4").Value = "2016-03-02 14:52:36"
$wsZhCn.Range("G4").Value = "2016-03-02 14:53:18"
$wsZhCn.Range("H4").Value = $reason

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4a5b6c7d8e9f0a1b2c3d4e5f60718293a4b5c6d7/e2e/$newMd",
    "",
    "",
    $newMd
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b6c7d8e9f0a1b2c3d4e5f60718293a4b5c6d7e8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "",
    "",
    $zhXlf
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4a5b6c7d8e9f0a1b2c3d4e5f60718293a4b5c6d7/e2e/$newMd",
    "",
    "",
    $newMd
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6c7d8e9f0a1b2c3d4e5f60718293a4b5c6d7e8f9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "",
    "",
    $zhXlf
)

$wsZhCn.Range("A4").Style = $wsZhCn.Range("A3").Style
$wsZhCn.Range("B4").Style = $wsZhCn.Range("B3").Style
$wsZhCn.Range("C4").Style = $wsZhCn.Range("C3").Style
$wsZhCn.Range("D4").Style = $wsZhCn.Range("D3").Style
$wsZhCn.Range("E4").Style = $wsZhCn.Range("E3").Style
$wsZhCn.Range("F4").Style = $wsZhCn.Range("F3").Style
$wsZhCn.Range("G4").Style = $wsZhCn.Range("G3").Style
$wsZhCn.Range("H4").Style = $wsZhCn.Range("H3").Style

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf = "$newXlfBase.de-de.xlf"

$wsDeDe.Range("B4").Value = $status
$wsDeDe.Range("D4").Value = "2016-03-02 14:52:45"
$wsDeDe.Range("G4").Value = "2016-03-02 14:53:39"
$wsDeDe.Range("H4").Value = $reason

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7d8e9f0a1b2c3d4e5f60718293a4b5c6d7e8f9a0/e2e/$newMd",
    "",
    "",
    $newMd
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e9f0a1b2c3d4e5f60718293a4b5c6d7e8f9a0b1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "",
    "",
    $deXlf
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7d8e9f0a1b2c3d4e5f60718293a4b5c6d7e8f9a0/e2e/$newMd",
    "",
    "",
    $newMd
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9f0a1b2c3d4e5f60718293a4b5c6d7e8f9a0b1c2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "",
    "",
    $deXlf
)

$wsDeDe.Range("A4").Style = $wsDeDe.Range("A3").Style
$wsDeDe.Range("B4").Style = $wsDeDe.Range("B3").Style
$wsDeDe.Range("C4").Style = $wsDeDe.Range("C3").Style
$wsDeDe.Range("D4").Style = $wsDeDe.Range("D3").Style
$wsDeDe.Range("E4").Style = $wsDeDe.Range("E3").Style
$wsDeDe.Range("F4").Style = $wsDeDe.Range("F3").Style
$wsDeDe.Range("G4").Style = $wsDeDe.Range("G3").Style
$wsDeDe.Range("H4").Style = $wsDeDe.Range("H3").Style

Write-Output "Added handback row for 185fe346-231a-4d9d-9928-1b9795caeec0"
